$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: force Text format on all Price (D) cells being updated,
# so numeric-looking strings (e.g. "1.003", "90.90") are preserved as literal text
# instead of being parsed/rounded as numbers by Excel.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Step 2: write the updated cell values for each changed row
# Row 2
$ws.Range("D2").Value = "22.464.16"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$ws.Range("D3").Value = "1.571.23"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "1.005"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6
$ws.Range("D6").Value = "290.84"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").Value = "0.3752"
$ws.Range("E7").Value = "  -0.23%  "

# Row 8
$ws.Range("D8").Value = "49.62"
$ws.Range("E8").Value = "  -0.94%  "

# Row 9
$ws.Range("D9").Value = "0.3407"
$ws.Range("E9").Value = "  -0.28%  "

# Row 10
$ws.Range("D10").Value = "1.146"
$ws.Range("E10").Value = "  -2.03%  "

# Row 11
$ws.Range("D11").Value = "0.07549"
$ws.Range("E11").Value = "  -1.71%  "

# Row 12
$ws.Range("D12").Value = "1.007"
$ws.Range("E12").Value = "  +0.52%  "

# Row 13
$ws.Range("D13").Value = "21.32"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14
$ws.Range("D14").Value = "5.982"
$ws.Range("E14").Value = "  -0.26%  "

# Row 15
$ws.Range("D15").Value = "6.931"
$ws.Range("E15").Value = "  -0.16%  "

# Row 16
$ws.Range("D16").Value = "1.583.79"
$ws.Range("E16").Value = "  +0.64%  "

# Row 17
$ws.Range("D17").Value = "0.00001117"
$ws.Range("E17").Value = "  -2.54%  "

# Row 18
$ws.Range("D18").Value = "90.90"
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$ws.Range("D19").Value = "0.06768"
$ws.Range("E19").Value = "  +0.59%  "

# Row 20
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").Value = "6.192"
$ws.Range("E21").Value = "  -0.82%  "

# Row 22
$ws.Range("D22").Value = "16.33"
$ws.Range("E22").Value = "  -2.61%  "

# Row 23
$ws.Range("D23").Value = "12.08"
$ws.Range("E23").Value = "  +0.43%  "

# Row 24
$ws.Range("D24").Value = "22.384.23"
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("D25").Value = "2.391"
$ws.Range("E25").Value = "  -0.14%  "

# Row 26
$ws.Range("D26").Value = "2.615"
$ws.Range("E26").Value = "  -6.18%  "

# Row 27
$ws.Range("D27").Value = "20.19"
$ws.Range("E27").Value = "  -0.45%  "

# Row 28
$ws.Range("D28").Value = "148.16"
$ws.Range("E28").Value = "  +2.41%  "

# Row 29
$ws.Range("D29").Value = "5.002"
$ws.Range("E29").Value = "  -1.15%  "

# Row 30
$ws.Range("D30").Value = "126.00"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("D31").Value = "1.757.14"
$ws.Range("E31").Value = "  +0.59%  "

# Row 32
$ws.Range("D32").Value = "1.035"
$ws.Range("E32").Value = "  -0.21%  "

# Row 33
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "1.994"
$ws.Range("E33").Value = "  -1.50%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "6.076"
$ws.Range("E34").Value = "  -3.27%  "

# Row 35
$ws.Range("D35").Value = "9.859"
$ws.Range("E35").Value = "  -2.90%  "

# Row 36
$ws.Range("D36").Value = "0.08419"
$ws.Range("E36").Value = "  -1.44%  "

# Row 37
$ws.Range("D37").Value = "1.371"
$ws.Range("E37").Value = "  +5.54%  "

# Row 38
$ws.Range("D38").Value = "0.02453"
$ws.Range("E38").Value = "  -3.92%  "

# Row 39
$ws.Range("D39").Value = "0.06570"
$ws.Range("E39").Value = "  +0.18%  "

# Row 40
$ws.Range("D40").Value = "0.2270"
$ws.Range("E40").Value = "  -2.63%  "

# Row 41
$ws.Range("D41").Value = "5.453"
$ws.Range("E41").Value = "  -1.26%  "

# Row 42
$ws.Range("D42").Value = "11.38"
$ws.Range("E42").Value = "  -2.73%  "

# Row 43
$ws.Range("D43").Value = "0.6264"
$ws.Range("E43").Value = "  -2.94%  "

# Row 44
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.21%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "14.04"
$ws.Range("E45").Value = "  -0.49%  "

# Row 46
$ws.Range("D46").Value = "3.815"
$ws.Range("E46").Value = "  +0.92%  "

# Row 47
$ws.Range("D47").Value = "0.5863"
$ws.Range("E47").Value = "  -2.99%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "130.29"
$ws.Range("E48").Value = "  +3.75%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.093"
$ws.Range("E49").Value = "  -0.53%  "

# Row 50
$ws.Range("D50").Value = "1.221"
$ws.Range("E50").Value = "  -6.56%  "

# Row 51
$ws.Range("D51").Value = "0.07315"
$ws.Range("E51").Value = "  -0.21%  "

# Step 3: reset style index on Price (D) cells back to the default/Normal style
# (clears the temporary Text number-format flag without altering visible style)
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}